# Code update for Origin Field
# - Copies row 9 of the "View Load" sheet into a new row 10 (same template
#   row used by the automation framework), then updates the Origin (G) and
#   Status (F) values for the two rows to the latest generated Origin codes.
# - Makes "View Load" the active/selected sheet/tab, with the cursor resting
#   on E7 (matching the state captured after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("View Load")

# Update the Origin value on the existing row (row 9).
$ws.Range("G9").Value = "Alaska_1006061652"

# Duplicate row 9 into row 10, preserving formatting (number format on the
# Rate column, etc.), then adjust the two cells that differ for the new row.
$ws.Range("A9:N9").Copy($ws.Range("A10:N10"))
$ws.Range("F10").Value = "Open"
$ws.Range("G10").Value = "Alaska_1006055219"

# Make "View Load" the active sheet/tab and move the selection to E7.
$ws.Activate()
$ws.Range("E7").Select()
